# Apply the committed change to the "params" sheet:
#  - clear the contents of row 7 (the "carbon_intensity" / "S1" scenario override row)
#  - add a new "id" column (S) with sequential numbers 1..5 for the remaining data rows (2..6)
#  - move the sheet selection to S10 (just below the new data)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Clear contents of row 7 (A7:O7) but keep existing formatting/styles.
$ws.Range("A7:O7").ClearContents()

# Populate the new "id" column (S) for rows 2-6.
$ws.Range("S2").Value = 1
$ws.Range("S3").Value = 2
$ws.Range("S4").Value = 3
$ws.Range("S5").Value = 4
$ws.Range("S6").Value = 5

# Update the active selection on the sheet.
$ws.Range("S10").Select()
